$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Caso1")
$ws.Range("G2").Value = 0.9687058329582214
$ws.Range("G3").Value = 0.969495415687561
$ws.Range("G4").Value = 0.9671235084533691
$ws.Range("G5").Value = 0.9749326705932617
$ws.Range("G6").Value = 0.9757153987884521
$ws.Range("G7").Value = 0.9670174717903137
$ws.Range("G8").Value = 0.9654960632324219
$ws.Range("G9").Value = 0.9679502844810486
$ws.Range("G10").Value = 0.9691173434257507
$ws.Range("G11").Value = 0.9680575132369995
$ws.Range("G12").Value = 0.9657676815986633
$ws.Range("G13").Value = 0.9655358195304871
$ws.Range("G14").Value = 0.9699353575706482
$ws.Range("G15").Value = 0.9740872979164124
$ws.Range("G16").Value = 0.974168598651886
$ws.Range("G17").Value = 0.9630635976791382
$ws.Range("G18").Value = 0.9666027426719666
$ws.Range("G19").Value = 0.968299925327301

$ws = $wb.Worksheets.Item("Caso2")
$ws.Range("G2").Value = 0.9697662591934204
$ws.Range("G3").Value = 0.9701190590858459
$ws.Range("G4").Value = 0.9674259424209595
$ws.Range("G5").Value = 0.9759608507156372
$ws.Range("G6").Value = 0.9761242866516113
$ws.Range("G7").Value = 0.9669963717460632
$ws.Range("G8").Value = 0.9664857387542725
$ws.Range("G9").Value = 0.968216598033905
$ws.Range("G10").Value = 0.969961404800415
$ws.Range("G11").Value = 0.9685810208320618
$ws.Range("G12").Value = 0.9658108949661255
$ws.Range("G13").Value = 0.9662909507751465
$ws.Range("G14").Value = 0.970343291759491
$ws.Range("G15").Value = 0.9747560024261475
$ws.Range("G16").Value = 0.9753343462944031
$ws.Range("G17").Value = 0.9634864926338196
$ws.Range("G18").Value = 0.9672340154647827
$ws.Range("G19").Value = 0.9688071012496948

$ws = $wb.Worksheets.Item("Caso3")
$ws.Range("G2").Value = 0.9690381288528442
$ws.Range("G3").Value = 0.9698264598846436
$ws.Range("G4").Value = 0.9672456979751587
$ws.Range("G5").Value = 0.9754802584648132
$ws.Range("G6").Value = 0.9756767749786377
$ws.Range("G7").Value = 0.9664040803909302
$ws.Range("G8").Value = 0.9657930731773376
$ws.Range("G9").Value = 0.9675602912902832
$ws.Range("G10").Value = 0.9694263339042664
$ws.Range("G11").Value = 0.96832275390625
$ws.Range("G12").Value = 0.9656330943107605
$ws.Range("G13").Value = 0.9653477668762207
$ws.Range("G14").Value = 0.9700199961662292
$ws.Range("G15").Value = 0.9744585752487183
$ws.Range("G16").Value = 0.9746558666229248
$ws.Range("G17").Value = 0.9631686210632324
$ws.Range("G18").Value = 0.9673013687133789
$ws.Range("G19").Value = 0.9682842493057251
